# Actualización automática del mapa (2025-08-15 11:42:08)
# Applies the changes described by the upstream diff:
#  - Row 5: updated case number / address / observation / attachments flag,
#           coordinate & zone columns cleared (no longer geocoded).
#  - Row 72: updated case number / address / observation,
#           coordinate & zone columns cleared (no longer geocoded).
#  - Three brand-new pending rows (88-90) appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to stay TEXT even when
# the content looks numeric / date-like (Excel would otherwise silently
# re-type "6968" as a Number or "8/14/2025" as a date serial). We flip the
# cell to text format, assign the literal value, then strip the formatting
# change back off so no stray style index is left behind on the cell.
function Set-TextValue {
    param($cell, $value)
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# ---------------------------------------------------------------------
# Row 5
# ---------------------------------------------------------------------
Set-TextValue "A5" "6968"
Set-TextValue "C5" "SANTA FE AV. 5154"
$ws.Range("H5").Value = "Picada"
$ws.Range("I5").Value = 1
$ws.Range("M5:P5").ClearContents()

# ---------------------------------------------------------------------
# Row 72
# ---------------------------------------------------------------------
Set-TextValue "A72" "6942"
Set-TextValue "C72" "3 DE FEBRERO 2169"
$ws.Range("H72").Value = "Picada"
$ws.Range("M72:P72").ClearContents()

# ---------------------------------------------------------------------
# New rows 88-90
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row = 88; A = "6947"; B = "8/14/2025"; C = "ALMAFUERTE AV. 682";      D = "4"; E = "808972978"; H = "Cambiar" },
    @{ Row = 89; A = "6951"; B = "8/14/2025"; C = "MEXICO 2751";             D = "3"; E = "808972984"; H = "Picada"  },
    @{ Row = 90; A = "6960"; B = "8/14/2025"; C = "VALLESE, FELIPE 1940";    D = "7"; E = "808972988"; H = "Picada"  }
)

foreach ($r in $newRows) {
    $row = $r.Row
    Set-TextValue "A$row" $r.A
    Set-TextValue "B$row" $r.B
    Set-TextValue "C$row" $r.C
    Set-TextValue "D$row" $r.D
    Set-TextValue "E$row" $r.E
    $ws.Range("F$row").Value = "AYKO"
    $ws.Range("G$row").Value = "Pendiente"
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = 1
    $ws.Range("J$row").Value = "Cambio"
    $ws.Range("K$row").Value = "Sin equipos"
    $ws.Range("L$row").Value = "Pasante"
    $ws.Range("M$row:P$row").ClearContents()
}
